$d = $word.ActiveDocument

# --- 1. Insert a new "Meta description" paragraph right after the H1 title ---
$titlePara = $d.Paragraphs.Item(1)
$breakRange = $titlePara.Range.Duplicate
$breakRange.Collapse(0)
$breakRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$fullNewRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Beat the Beast Mighty Sphinx and play for free. Enjoy simple gameplay, high volatility, and impressive graphics.</w:t></w:r></w:p>'
$fullNewRange.InsertXML($metaXml)

# --- 2. Remove the duplicated bold title paragraph near the end of the document ---
$count = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($count - 1)
$secondToLast.Range.Delete()

# --- 3. Replace the final (italic) paragraph's text with the new DALLE prompt ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastRange = $lastPara.Range
$textRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$textRange.Text = "Prompt: DALLE, create a cartoon-style feature image for `"Beat the Beast Mighty Sphinx`" featuring a happy Maya warrior with glasses. Description: The feature image should be in cartoon-style with bright and vibrant colors. The main element of the image should be a happy and confident Maya warrior with glasses, standing in front of a giant Sphinx. The warrior should be wearing traditional Maya clothing, with a feather headdress and accessories. The background should have an Egyptian theme, with hieroglyphics and pyramids visible. The Sphinx should be portrayed as dark and imposing, with glowing yellow eyes. The image should emphasize the adventure, excitement, and mystery of Ancient Egypt, while also showcasing the unique blend of Maya and Egyptian elements in the game."

Write-Output "done"
